$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row updates: Price (column D) and Volume(1h) (column E) changes.
# D is $null when the price cell itself did not change in this update.
$updates = @(
    @{Row=2; D="67.042.50"; E="  +0.59%  "},
    @{Row=3; D="3.505.36"; E="  -0.04%  "},
    @{Row=5; D="594.48"; E="  +0.47%  "},
    @{Row=6; D="173.05"; E="  +1.63%  "},
    @{Row=8; D="0.600"; E="  +1.53%  "},
    @{Row=9; D=$null; E="  +4.14%  "},
    @{Row=10; D="7.28"; E="  -0.85%  "},
    @{Row=11; D=$null; E="  -1.36%  "},
    @{Row=12; D="4.111.53"; E="  -0.02%  "},
    @{Row=13; D=$null; E="  -0.16%  "},
    @{Row=14; D="29.06"; E="  +2.22%  "},
    @{Row=15; D="67.005.10"; E="  +0.48%  "},
    @{Row=16; D="0.0000179"; E="  +0.22%  "},
    @{Row=17; D="3.492.26"; E="  -0.52%  "},
    @{Row=18; D=$null; E="  -0.60%  "},
    @{Row=19; D=$null; E="  +0.27%  "},
    @{Row=20; D="394.60"; E="  +0.98%  "},
    @{Row=21; D="8.02"; E="  +0.18%  "},
    @{Row=22; D="73.13"; E="  +0.00%  "},
    @{Row=23; D="0.999"; E="  -0.04%  "},
    @{Row=24; D=$null; E="  -0.25%  "},
    @{Row=25; D=$null; E="  -2.89%  "},
    @{Row=26; D=$null; E="  -1.46%  "},
    @{Row=27; D="10.23"; E="  -1.14%  "},
    @{Row=28; D=$null; E="  +0.61%  "},
    @{Row=29; D=$null; E="  -0.25%  "},
    @{Row=30; D="6.30"; E="  -1.31%  "},
    @{Row=31; D=$null; E="  -3.60%  "},
    @{Row=32; D=$null; E="  -0.45%  "},
    @{Row=33; D="23.79"; E="  +0.61%  "},
    @{Row=34; D="7.36"; E="  -1.11%  "},
    @{Row=35; D="1.67"; E="  +2.88%  "},
    @{Row=36; D="163.55"; E="  +0.40%  "},
    @{Row=37; D="0.882"; E="  -0.29%  "},
    @{Row=38; D=$null; E="  -0.26%  "},
    @{Row=39; D="7.02"; E="  +3.36%  "},
    @{Row=40; D="4.68"; E="  -0.88%  "},
    @{Row=41; D="0.0747"; E="  -0.17%  "},
    @{Row=42; D="27.50"; E="  +0.88%  "},
    @{Row=43; D="26.36"; E="  -0.89%  "},
    @{Row=44; D="2.812.25"; E="  -0.11%  "},
    @{Row=45; D="2.57"; E="  +1.58%  "},
    @{Row=46; D="42.71"; E="  -0.97%  "},
    @{Row=47; D=$null; E="  -2.67%  "},
    @{Row=48; D="338.17"; E="  -4.90%  "},
    @{Row=49; D="34.71"; E="  +2.77%  "},
    @{Row=50; D=$null; E="  -0.85%  "},
    @{Row=51; D=$null; E="  -1.07%  "}
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        # Force the Price column to remain plain text (it holds values like
        # "67.042.50" that are not valid numbers, as well as values like
        # "0.600" whose trailing zero must be preserved as text). Setting a
        # text number-format before assigning the value prevents Excel from
        # auto-converting the string into a numeric value, and resetting the
        # style back to Normal afterwards keeps the cell's style index
        # unchanged (matching the original workbook, which carries no
        # explicit style on these cells).
        $cellD = $ws.Cells.Item($row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
        $cellD.Style = "Normal"
    }

    $ws.Cells.Item($row, 5).Value = $u.E
}
